$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 08.06.2025"

$ws.Range("B6").Value = "12.06."
$ws.Range("C6").Value = "13.06."
$ws.Range("D6").Value = "AMAZON.DE MKTPLC EU SYXWKD"
$ws.Range("E6").Value = "107,98-"

$ws.Range("B7").Value = "16.06."
$ws.Range("C7").Value = "17.06."
$ws.Range("D7").Value = "ZALANDO MKTPLC EU CHGFLY"
$ws.Range("E7").Value = "115,90-"

$ws.Range("B8").Value = "17.06."
$ws.Range("C8").Value = "18.06."
$ws.Range("E8").Value = "25,32-"

$ws.Range("D12").Value = "KONTOSTAND AM 22.06.2025"
$ws.Range("E12").Value = "249,20-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 29.06.2025"
